$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.083.02"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.755.73"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'337.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.3773"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.25%  "
$ws.Range("D8").Value = "'0.3352"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.89%  "
$ws.Range("D9").Value = "'45.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.93%  "
$ws.Range("D10").Value = "'1.116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.37%  "
$ws.Range("E11").Value = "  -4.57%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'22.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E14").Value = "  -5.68%  "
$ws.Range("D15").Value = "'7.137"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "1.758.09"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "'0.00001056"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").Value = "'0.06585"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "'80.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.68%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").Value = "'6.234"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("D23").Value = "28.079.80"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  -5.42%  "
$ws.Range("D25").Value = "'2.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").Value = "'152.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "'19.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.63%  "
$ws.Range("D28").Value = "'2.318"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.52%  "
$ws.Range("D29").Value = "1.956.77"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "'131.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("E31").Value = "  -15.06%  "
$ws.Range("D32").Value = "'4.023"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "'5.773"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.77%  "
$ws.Range("D34").Value = "'0.08725"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("D35").Value = "'12.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.68%  "
$ws.Range("D36").Value = "'0.6664"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("D37").Value = "'0.02329"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D40").Value = "'0.2109"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").Value = "'1.214"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "'1.443"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.25%  "
$ws.Range("D43").Value = "'8.022"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.66%  "
$ws.Range("D45").Value = "'13.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.26%  "
$ws.Range("D46").Value = "'3.837"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'0.6043"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.47%  "
$ws.Range("D48").Value = "'128.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "'2.014"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.92%  "

# Row 38/39 swap (InternetComputer(DFINITY) <-> Hedera), with updated D/E values
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.157"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.92%  "

# Row 50/51 swap (EOS <-> Cronos), with updated D/E values
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07147"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.175"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "

Write-Output "cryptos list updated"
